# Updated cryptos list on Wed Jul 31 06:29:53 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.404.95"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.317.16"
$ws.Range("E3").Value = "  -0.04%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB (D looks numeric - force text like the source cell)
$ws.Range("D5").Value = "'586.71"
$ws.Range("E5").Value = "  +2.43%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'182.64"
$ws.Range("E6").Value = "  +0.25%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.643"
$ws.Range("E7").Value = "  +7.60%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -2.12%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +2.41%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +0.18%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "3.893.29"
$ws.Range("E12").Value = "  -0.06%  "

# Row 14 - WrappedBTC
$ws.Range("D14").Value = "66.423.23"
$ws.Range("E14").Value = "  -0.40%  "

# Row 15 - Avalanche
$ws.Range("E15").Value = "  -2.69%  "

# Row 16 & 17 swap: ShibaInu <-> WrappedEther with new values
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.339.22"
$ws.Range("E16").Value = "  +0.73%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000164"
$ws.Range("E17").Value = "  -1.95%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "'431.58"
$ws.Range("E18").Value = "  -0.25%  "

# Row 19 - Chainlink
$ws.Range("D19").Value = "'13.33"
$ws.Range("E19").Value = "  -2.33%  "

# Row 20 - Polkadot
$ws.Range("E20").Value = "  -2.56%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "'7.43"
$ws.Range("E21").Value = "  -2.55%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "'72.33"
$ws.Range("E22").Value = "  -1.52%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.06%  "

# Row 24 - LEO
$ws.Range("D24").Value = "'5.69"
$ws.Range("E24").Value = "  +0.50%  "

# Row 25 - WrappedeETH
$ws.Range("D25").Value = "3.436.17"
$ws.Range("E25").Value = "  -0.73%  "

# Row 26 - Polygon
$ws.Range("D26").Value = "'0.517"
$ws.Range("E26").Value = "  -0.80%  "

# Row 27 - Kaspa
$ws.Range("E27").Value = "  +2.57%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  -3.49%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("E29").Value = "  -0.38%  "

# Row 30 - Binance-PegBSC-USD
$ws.Range("E30").Value = "  -0.08%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "'1.95"
$ws.Range("E31").Value = "  -0.52%  "

# Row 32 - EthereumClassic
$ws.Range("E32").Value = "  -1.47%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  +0.03%  "

# Row 34 - NEARProtocol
$ws.Range("E34").Value = "  -1.53%  "

# Row 35 & 36 swap: Fetch.AI <-> Aptos with new values
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "'6.65"
$ws.Range("E35").Value = "  -2.33%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "'1.21"
$ws.Range("E36").Value = "  -2.82%  "

# Row 37 & 38 swap: Monero <-> ImmutableX with new values
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.47"
$ws.Range("E37").Value = "  -2.73%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'159.08"
$ws.Range("E38").Value = "  -0.56%  "

# Row 39 - Stacks
$ws.Range("D39").Value = "'1.82"
$ws.Range("E39").Value = "  -1.19%  "

# Row 40 - EnergySwap
$ws.Range("E40").Value = "  -1.28%  "

# Row 41 - Maker
$ws.Range("D41").Value = "2.873.40"
$ws.Range("E41").Value = "  +1.92%  "

# Row 42 - Mantle
$ws.Range("D42").Value = "'0.773"
$ws.Range("E42").Value = "  -2.19%  "

# Row 43 - Filecoin
$ws.Range("E43").Value = "  -2.00%  "

# Row 44 - OKB
$ws.Range("D44").Value = "'40.24"
$ws.Range("E44").Value = "  +0.19%  "

# Row 45 - Hedera
$ws.Range("E45").Value = "  -1.15%  "

# Row 46 - RenderToken
$ws.Range("D46").Value = "'6.03"
$ws.Range("E46").Value = "  -2.48%  "

# Row 47 - dogwifhat
$ws.Range("E47").Value = "  -1.58%  "

# Row 48 & 49 swap: InjectiveProtocol <-> Bittensor with new values
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "'319.03"
$ws.Range("E48").Value = "  -1.76%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "'23.47"
$ws.Range("E49").Value = "  -3.76%  "

# Row 50 - VeChain
$ws.Range("E50").Value = "  -0.49%  "

# Row 51 - Stellar
$ws.Range("E51").Value = "  +3.69%  "
